$wb = $excel.ActiveWorkbook

# --- 1. Status text update: "Ready for handoff" -> "In Translation" ------
# This shows up in several cells (Overview!E2/F2 per-language status, and the
# Status column on each language sheet). Walk every sheet's used range and
# replace the exact (string-typed) match so a boolean "True"/"False" cell
# elsewhere can't accidentally satisfy a loose comparison.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if (($val -is [string]) -and ($val -eq "Ready for handoff")) {
                $cell.Value = "In Translation"
            }
        }
    }
}

# --- 2. Re-fit the now-narrower status columns ----------------------------
# With shorter text ("In Translation" vs "Ready for handoff") the status
# columns get re-sized: Overview columns E & F, and column C ("Status") on
# both the zh-cn and de-de sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
